$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) data for rows 2..29: row#, date-serial, country, confirmed, new, deaths, new-deaths.
# This is the fully re-sorted table (by date, then original country order) with the
# newly-added NL and France figures folded in, plus one late-appended France row (r29)
# that was not swept up by the A2:F28 sort range.
$rows = @(
    @(2, 43894, "Italy", 2502, 466, 80, 28),
    @(3, 43894, "UK", 51, 12, 0, 0),
    @(4, 43894, "NL", 28, 10, 0, 0),
    @(5, 43895, "Italy", 3089, 587, 107, 27),
    @(6, 43895, "UK", 89, 38, 0, 0),
    @(7, 43895, "NL", 38, 10, 0, 0),
    @(8, 43895, "France", 282, 73, 4, 0),
    @(9, 43896, "Italy", 3858, 769, 148, 41),
    @(10, 43896, "UK", 118, 29, 0, 0),
    @(11, 43896, "NL", 82, 44, 0, 0),
    @(12, 43896, "France", 420, 138, 6, 2),
    @(13, 43897, "Italy", 4636, 778, 197, 49),
    @(14, 43897, "UK", 167, 49, 1, 1),
    @(15, 43897, "NL", 128, 46, 1, 1),
    @(16, 43897, "France", 613, 193, 9, 3),
    @(17, 43898, "Italy", 5883, 1247, 234, 37),
    @(18, 43898, "UK", 210, 43, 2, 1),
    @(19, 43898, "NL", 188, 60, 1, 0),
    @(20, 43898, "France", 706, 93, 10, 1),
    @(21, 43899, "Italy", 7375, 1492, 366, 132),
    @(22, 43899, "UK", 277, 67, 2, 0),
    @(23, 43899, "NL", 265, 77, 3, 2),
    @(24, 43899, "France", 1116, 410, 19, 9),
    @(25, 43900, "Italy", 9172, 1797, 463, 97),
    @(26, 43900, "UK", 323, 46, 3, 1),
    @(27, 43900, "NL", 321, 56, 3, 0),
    @(28, 43900, "France", 1402, 286, 30, 11),
    @(29, 43894, "France", 212, 21, 4, 1)
)

# Column A already carries the short-date style (from the pre-existing rows) - extend
# that same number format to the brand-new rows below the original A15 by copying the
# format (not the value) from the first data cell down across the whole column range.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A2:A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Re-run the Data > Sort (ascending by date) over the same A2:F28 window the author
# used - the new rows above are already in final order, so this is a no-op on the
# values, but it is what records the worksheet's remembered sort range/state.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1"))
$sortObj.SetRange($ws.Range("A2:F28"))
$sortObj.Header = -4142
$sortObj.Apply()

# Match the saved page setup / selection state from the edit.
$ws.PageSetup.Orientation = 1
$ws.Range("F29").Select()
